# Update the "取得日時" (retrieved datetime) column on the "ランサーズ" sheet
# for the already-listed rows (2-13) to the new timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-01-16 12:53:56"

for ($row = 2; $row -le 13; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
